$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cells($ws, $values) {
    foreach ($key in $values.Keys) {
        $val = $values[$key]
        if ($null -eq $val) {
            $ws.Range($key).Value = ""
        } else {
            $ws.Range($key).Value = $val
        }
    }
}

# Row 2
$row2 = @{
    "D2" = 0.2225
    "E2" = 0.239
    "F2" = 0.145
    "I2" = 0.000552023048734038
    "J2" = 0.0003710785467716233
    "K2" = 1737.1
    "L2" = 0.4580114430353047
    "M2" = 599.92
    "N2" = 0.06079017499772005
    "O2" = 0.3453572045362961
    "P2" = 599.92
    "Q2" = 0.06079017499772005
    "R2" = 0.3453572045362961
    "U2" = 7484.599999999999
    "V2" = 0.7584180287170548
    "W2" = 0.2343829373561452
    "X2" = 0.06398079505317383
    "Y2" = 0.1704021423029713
    "Z2" = 2.416561024472552
    "AB2" = 0.06262188397835444
    "AC2" = -0.06262188397835444
    "AD2" = 1319.8
    "AE2" = 55.53171091533207
    "AF2" = 1375.331710915332
    "AG2" = -6109.268289084667
    "AH2" = 0.1223165983763818
    "AI2" = 0.1341861431644689
    "AJ2" = -1.625051007402713
    "AK2" = -2.209634772693675
    "AN2" = 99.98484848484848
    "AP2" = -462.8233552336869
}
Set-Cells $ws $row2

# Row 3
$row3 = @{
    "D3" = 0.166
    "E3" = 0.244
    "K3" = 135.7
    "L3" = 0.4718358831710708
    "M3" = 47.8
    "N3" = 0.1062694530902623
    "O3" = 0.3522476050110538
    "P3" = 47.8
    "Q3" = 0.1062694530902623
    "R3" = 0.3522476050110538
    "U3" = 1289.0
    "V3" = 2.86571809693197
    "W3" = 0.1681953396132871
    "X3" = 0.05829854405455182
    "Y3" = 0.1098967955587352
    "Z3" = -0.1913926544084862
    "AB3" = 0.05817395895350823
    "AC3" = -0.05817395895350823
    "AD3" = 10.0
    "AF3" = 10.0
    "AG3" = -1279.0
    "AH3" = 0.02174858634188778
    "AI3" = 0.01077702338614075
    "AJ3" = 1.542450554751568
    "AK3" = 3.541955137081141
}
Set-Cells $ws $row3

# Row 4
$row4 = @{
    "B4" = "Credit Agricole - Egypt Bank (S.A.E.) (CASE:CIEB)"
    "D4" = 0.0971
    "E4" = 0.0966
    "F4" = 0.0613
    "I4" = 0.0
    "J4" = 0.0
    "K4" = 97.3
    "L4" = 0.450462962962963
    "M4" = 81.7
    "N4" = 0.142533147243545
    "O4" = 0.8396711202466599
    "P4" = 81.7
    "Q4" = 0.142533147243545
    "R4" = 0.8396711202466599
    "U4" = 825.4
    "V4" = 1.439986043265876
    "W4" = 0.2400690846286701
    "X4" = 0.0599890110554402
    "Y4" = 0.1800800735732299
    "Z4" = -1.30988477865373
    "AA4" = -0.0
    "AB4" = 0.05951979869893757
    "AC4" = -0.05951979869893757
    "AD4" = 38.7
    "AE4" = 0.0
    "AF4" = 38.7
    "AG4" = -786.6999999999999
    "AH4" = 0.06324562837064879
    "AI4" = 0.0817490494296578
    "AJ4" = 3.684777517564405
    "AK4" = 2.234943181818182
    "AN4" = $null
    "AP4" = $null
}
Set-Cells $ws $row4

# Row 5
$row5 = @{
    "B5" = "QNB ALAHLI Bank (S.A.E) (CASE:QNBA)"
    "D5" = 0.224
    "E5" = 0.212
    "F5" = $null
    "I5" = 0.001689049262589023
    "J5" = 0.001202424063773522
    "K5" = 489.0
    "L5" = 0.5231065468549423
    "M5" = 192.4
    "N5" = 0.0781700727257953
    "O5" = 0.3934560327198364
    "P5" = 192.4
    "Q5" = 0.0781700727257953
    "R5" = 0.3934560327198364
    "U5" = 862.5
    "V5" = 0.3504245723804493
    "W5" = 0.2390730419477853
    "X5" = 0.06362347673016834
    "Y5" = 0.1754495652176169
    "Z5" = 0.6455331293509919
    "AA5" = 0.0007762045686946582
    "AB5" = 0.06205900294237533
    "AC5" = -0.06128279837368067
    "AD5" = 384.1
    "AE5" = 21.70538374665891
    "AF5" = 405.8053837466589
    "AG5" = -456.6946162533411
    "AH5" = 0.1415383564368196
    "AI5" = 0.1443014745978512
    "AJ5" = -0.2278227026407397
    "AK5" = -0.2342377571814117
    "AN5" = 64.88175675675676
    "AP5" = -77.14436085360491
}
Set-Cells $ws $row5

# Row 6
$row6 = @{
    "B6" = "Housing and Development Bank- Egypt (S.A.E) (CASE:HDBK)"
    "D6" = 0.23
    "E6" = 0.348
    "F6" = $null
    "K6" = 136.7
    "L6" = 0.4896131805157593
    "M6" = 34.3
    "N6" = 0.1092008914358484
    "O6" = 0.2509144111192392
    "P6" = 34.3
    "Q6" = 0.1092008914358484
    "R6" = 0.2509144111192392
    "U6" = 265.4
    "V6" = 0.8449538363578477
    "W6" = 0.3204406938584153
    "X6" = 0.06433811337617934
    "Y6" = 0.256102580482236
    "Z6" = 3.028199566160519
    "AB6" = 0.06250915410468716
    "AC6" = -0.06250915410468716
    "AD6" = 57.8
    "AF6" = 57.8
    "AG6" = -207.6
    "AH6" = 0.1554181231513848
    "AI6" = 0.0950814278664254
    "AJ6" = -1.949295774647886
    "AK6" = -0.6061313868613136
}
Set-Cells $ws $row6

# Row 7
$row7 = @{
    "B7" = "Commercial International Bank -Egypt S.A.E (CASE:COMI)"
    "D7" = 0.23
    "E7" = 0.204
    "F7" = 0.17
    "K7" = 673.0
    "L7" = 0.4363895733367916
    "M7" = 213.8
    "N7" = 0.03847884383492612
    "O7" = 0.3176820208023775
    "P7" = 213.8
    "Q7" = 0.03847884383492612
    "R7" = 0.3176820208023775
    "U7" = 3408.1
    "V7" = 0.613375807641776
    "W7" = 0.2296928327645051
    "X7" = 0.06090786639154877
    "Y7" = 0.1687849663729563
    "Z7" = 0.6994104308390023
    "AB7" = 0.06273461385202174
    "AC7" = -0.06273461385202174
    "AD7" = 511.9
    "AF7" = 511.9
    "AG7" = -2896.2
    "AH7" = 0.08435779967700471
    "AI7" = 0.1289160874382996
    "AJ7" = -1.088756061802188
    "AK7" = -5.146969966234225
}
Set-Cells $ws $row7

# Row 8
$row8 = @{
    "B8" = "Suez Canal Bank (S.A.E) (CASE:CANA)"
    "D8" = 0.262
    "E8" = $null
    "F8" = $null
    "I8" = 0.0008043182988565342
    "J8" = 0.0004208352709187439
    "K8" = 33.9
    "L8" = 0.3048561151079137
    "M8" = 3.06
    "N8" = 0.02326996197718631
    "O8" = 0.09026548672566372
    "P8" = 3.06
    "Q8" = 0.02326996197718631
    "R8" = 0.09026548672566372
    "U8" = 519.7
    "V8" = 3.952091254752852
    "W8" = 0.1751937984496124
    "X8" = 0.0740481906456032
    "Y8" = 0.1011456078040092
    "Z8" = -0.1653784397658025
    "AA8" = -0.00006959708050296064
    "AB8" = 0.06901588429456708
    "AC8" = -0.06908548137507003
    "AD8" = 53.7
    "AE8" = 4.702799025835767
    "AF8" = 58.40279902583577
    "AG8" = -461.2972009741643
    "AH8" = 0.3075404855822595
    "AI8" = 0.2043464907443272
    "AJ8" = 1.398729884946178
    "AK8" = 1.97222198065174
    "AN8" = 52.13592233009709
    "AP8" = -447.8613601690915
}
Set-Cells $ws $row8

# Row 9
$row9 = @{
    "B9" = "Abu Dhabi Islamic Bank - Egypt S.A.E. (CASE:ADIB)"
    "D9" = 0.208
    "E9" = 0.434
    "F9" = 0.145
    "I9" = 0.003191142466817553
    "J9" = 0.002173946422972891
    "K9" = 70.0
    "L9" = 0.3421309872922776
    "M9" = 1.16
    "N9" = 0.007641633728590249
    "O9" = 0.01657142857142857
    "P9" = 1.16
    "Q9" = 0.007641633728590249
    "R9" = 0.01657142857142857
    "T9" = 0.0
    "U9" = 51.8
    "V9" = 0.3412384716732542
    "W9" = 0.28
    "X9" = 0.08763942367521646
    "Y9" = 0.1923605763247836
    "Z9" = 0.7522460908566352
    "AA9" = 0.001635342698413123
    "AB9" = 0.07196466428491341
    "AC9" = -0.07032932158650029
    "AD9" = 104.1
    "AE9" = 18.58546125644564
    "AF9" = 122.6854612564456
    "AG9" = 70.88546125644564
    "AH9" = 0.4469652443333723
    "AI9" = 0.2720386171976451
    "AJ9" = 0.3183210114234337
    "AK9" = 0.1775752579598765
    "AN9" = 23.82151029748283
    "AP9" = 16.22092934930106
}
Set-Cells $ws $row9

# Row 10
$row10 = @{
    "D10" = 0.221
    "E10" = 0.239
    "I10" = -0.001048426426892449
    "J10" = -0.0007338984988247141
    "K10" = 101.5
    "L10" = 0.4675264854905574
    "M10" = 25.7
    "N10" = 0.1114000866926745
    "O10" = 0.2532019704433497
    "P10" = 25.7
    "Q10" = 0.1114000866926745
    "R10" = 0.2532019704433497
    "U10" = 262.7
    "V10" = 1.138708279150412
    "W10" = 0.2190332326283988
    "X10" = 0.084983276392565
    "Y10" = 0.1340499562358338
    "Z10" = -2.012758289537923
    "AA10" = 0.001477160287188881
    "AB10" = 0.07477489447326831
    "AC10" = -0.07329773418607943
    "AD10" = 159.5
    "AE10" = 10.53806688639175
    "AF10" = 170.0380668863918
    "AG10" = -92.66193311360823
    "AH10" = 0.4243122401810586
    "AI10" = 0.2360203807959049
    "AJ10" = -0.6712781133762975
    "AK10" = -0.202434404776316
    "AN10" = 84.8404255319149
    "AP10" = -49.28826229447247
}
Set-Cells $ws $row10

